$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 600, pushing the existing rows 600:641 down to 601:642
$ws.Rows.Item(600).EntireRow.Insert()

# Populate the newly inserted row 600 with the new day's first data point.
# The date column must stay plain text (matching every other row), so force
# a text format while assigning it, then drop that explicit formatting again
# so the cell ends up with the sheet's normal (unstyled) look.
$ws.Cells.Item(600, 1).NumberFormat = "@"
$ws.Cells.Item(600, 1).Value = "2026/01/07"
$ws.Cells.Item(600, 1).ClearFormats()

$ws.Cells.Item(600, 2).Value = "水"
$ws.Cells.Item(600, 3).Value = 17
$ws.Cells.Item(600, 4).Value = 201
